$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 3664.2144
$ws.Range("I2").Value = 657.9
$ws.Range("J2").Value = 5334.3887
$ws.Range("K2").Value = 657.9
$ws.Range("L2").Value = 5334.3887
$ws.Range("M2").Value = -544.9
$ws.Range("N2").Value = -5560.3887

$ws.Range("H9").Value = 131.63637
$ws.Range("I9").Value = 119.875
$ws.Range("K9").Value = 119.875
$ws.Range("M9").Value = 49.125

$ws.Range("H86").Value = 2922
$ws.Range("I86").Value = 2939.4285
$ws.Range("J86").Value = 2800
$ws.Range("K86").Value = 2939.4285
$ws.Range("L86").Value = 2800
$ws.Range("M86").Value = -1816.4285
$ws.Range("N86").Value = -5046

$ws.Range("H89").Value = 2922
$ws.Range("I89").Value = 2939.4285
$ws.Range("J89").Value = 2800
$ws.Range("K89").Value = 14697.1425
$ws.Range("L89").Value = 14000
$ws.Range("M89").Value = -9081.1425
$ws.Range("N89").Value = -25232

$ws.Range("H99").Value = 899.6667
$ws.Range("I99").Value = 899.6667
$ws.Range("K99").Value = 2699.0001
$ws.Range("M99").Value = -1201.0001

$ws.Range("H101").Value = 1095.25
$ws.Range("I101").Value = 555.6667
$ws.Range("K101").Value = 1667.0001
$ws.Range("M101").Value = -45.00009999999997

$ws.Range("H112").Value = 2369.4
$ws.Range("J112").Value = 2369.4
$ws.Range("L112").Value = 7108.200000000001
$ws.Range("N112").Value = -9324.200000000001

$ws.Range("H118").Value = 1153.3334
$ws.Range("I118").Value = 1153.3334
$ws.Range("J118").Value = 0
$ws.Range("K118").Value = 3460.0002
$ws.Range("L118").Value = 0
$ws.Range("M118").Value = -1803.0002
$ws.Range("N118").ClearContents()

$ws.Range("H138").Value = 1822.8387
$ws.Range("I138").Value = 1070.1875
$ws.Range("J138").Value = 2084.6304
$ws.Range("K138").Value = 3210.5625
$ws.Range("L138").Value = 6253.8912
$ws.Range("M138").Value = 1929.4375
$ws.Range("N138").Value = -16533.8912

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 2125.8147
$ws.Range("I63").Value = 1995.0952
$ws.Range("J63").Value = 2583.3333
$ws.Range("K63").Value = 1995.0952
$ws.Range("L63").Value = 2583.3333
$ws.Range("M63").Value = -1309.0952
$ws.Range("N63").Value = -3955.3333

$ws.Range("H66").Value = 2125.8147
$ws.Range("I66").Value = 1995.0952
$ws.Range("J66").Value = 2583.3333
$ws.Range("K66").Value = 9975.476000000001
$ws.Range("L66").Value = 12916.6665
$ws.Range("M66").Value = -6543.476000000001
$ws.Range("N66").Value = -19780.6665

$ws.Range("H122").Value = 10494.909
$ws.Range("J122").Value = 5170.8
$ws.Range("L122").Value = 15512.4
$ws.Range("N122").Value = -20412.4

$ws.Range("H132").Value = 4542.881
$ws.Range("I132").Value = 4911.8184
$ws.Range("J132").Value = 3190.111
$ws.Range("K132").Value = 14735.4552
$ws.Range("L132").Value = 9570.332999999999
$ws.Range("M132").Value = -12205.4552
$ws.Range("N132").Value = -14630.333

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2102.5
$ws.Range("I86").Value = 1996.3334
$ws.Range("J86").Value = 2341.375
$ws.Range("K86").Value = 1996.3334
$ws.Range("L86").Value = 2341.375
$ws.Range("M86").Value = -873.3334
$ws.Range("N86").Value = -4587.375

$ws.Range("H89").Value = 2102.5
$ws.Range("I89").Value = 1996.3334
$ws.Range("J89").Value = 2341.375
$ws.Range("K89").Value = 9981.666999999999
$ws.Range("L89").Value = 11706.875
$ws.Range("M89").Value = -4365.666999999999
$ws.Range("N89").Value = -22938.875

$ws.Range("H94").Value = 546.1
$ws.Range("I94").Value = 384.55554
$ws.Range("K94").Value = 384.55554
$ws.Range("M94").Value = 66.44445999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 495.625
$ws.Range("J107").Value = 535.2941
$ws.Range("L107").Value = 535.2941
$ws.Range("N107").Value = -4375.2941

$ws.Range("H122").Value = 1498.4546
$ws.Range("J122").Value = 1536.75
$ws.Range("L122").Value = 4610.25
$ws.Range("N122").Value = -9510.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 45623890
$ws.Range("I4").Value = 66914540
$ws.Range("K4").Value = 200743620
$ws.Range("M4").Value = -200743508

$ws.Range("H121").Value = 1428.92
$ws.Range("I121").Value = 303
$ws.Range("J121").Value = 1582.4546
$ws.Range("K121").Value = 909
$ws.Range("L121").Value = 4747.3638
$ws.Range("M121").Value = 401
$ws.Range("N121").Value = -7367.3638

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 116.04348
$ws.Range("I2").Value = 147.35715
$ws.Range("J2").Value = 67.333336
$ws.Range("K2").Value = 147.35715
$ws.Range("L2").Value = 67.333336
$ws.Range("M2").Value = -34.35714999999999
$ws.Range("N2").Value = -293.333336

$ws.Range("H11").Value = 35017500
$ws.Range("I11").Value = 33356666
$ws.Range("K11").Value = 33356666
$ws.Range("M11").Value = -33356527

$ws.Range("H80").Value = 3615.1667
$ws.Range("I80").Value = 3416.182
$ws.Range("J80").Value = 3927.8572
$ws.Range("K80").Value = 3416.182
$ws.Range("L80").Value = 3927.8572
$ws.Range("M80").Value = -2418.182
$ws.Range("N80").Value = -5923.8572

$ws.Range("H83").Value = 3615.1667
$ws.Range("I83").Value = 3416.182
$ws.Range("J83").Value = 3927.8572
$ws.Range("K83").Value = 17080.91
$ws.Range("L83").Value = 19639.286
$ws.Range("M83").Value = -12088.91
$ws.Range("N83").Value = -29623.286

$ws.Range("H122").Value = 2718.4614
$ws.Range("I122").Value = 2405.6667
$ws.Range("K122").Value = 7217.000100000001
$ws.Range("M122").Value = -4767.000100000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 9242
$ws.Range("I93").Value = 9025.429
$ws.Range("K93").Value = 9025.429
$ws.Range("M93").Value = -7777.429

$ws.Range("H122").Value = 46408.89
$ws.Range("I122").Value = 7565.2
$ws.Range("K122").Value = 22695.6
$ws.Range("M122").Value = -20245.6

$ws.Range("H132").Value = 5087.75
$ws.Range("I132").Value = 5016.1763
$ws.Range("K132").Value = 15048.5289
$ws.Range("M132").Value = -12518.5289

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 4755.467
$ws.Range("J81").Value = 6771.5
$ws.Range("L81").Value = 13543
$ws.Range("N81").Value = -15665

$ws.Range("H84").Value = 4755.467
$ws.Range("J84").Value = 6771.5
$ws.Range("L84").Value = 67715
$ws.Range("N84").Value = -78323

$ws.Range("H122").Value = 7385.4546
$ws.Range("J122").Value = 6949
$ws.Range("L122").Value = 20847
$ws.Range("N122").Value = -25747
